# "Generate Report for Handoff"
# Appends the handoff status of the newly discovered
# e398a89a-8f63-468c-a31f-1ce1b11d7198...md source file to the
# Overview / zh-cn / de-de localization-status report sheets.

$wb = $excel.ActiveWorkbook

$newBaseName   = "e398a89a-8f63-468c-a31f-1ce1b11d7198ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$newPathName   = "e2e\" + $newBaseName
$newHoUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/114f13da51ff14b176748d2c2f142bbf12d6a505/e2e/" + $newBaseName
$readyStatus   = "Ready for handoff"
$handoffDate   = "2016-08-29 20:39:14"
$zhXliffName   = "e398a89a-8f63-468c-a31f-1ce1b11d7198oooooooooooooooooooooooooooooooooooooooo.4e968e2ed0a8592eadc6eda2286ceddbd2510c41.zh-cn.xlf"
$deXliffName   = "e398a89a-8f63-468c-a31f-1ce1b11d7198oooooooooooooooooooooooooooooooooooooooo.4e968e2ed0a8592eadc6eda2286ceddbd2510c41.de-de.xlf"
$xliffDate     = "2016-08-29 20:39:09"
$epoch         = "0001-01-01 00:00:00"
$dateFmt       = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(3, 1).Value = $newBaseName
$wsOverview.Cells.Item(3, 2).Value = $newPathName
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 4).Value = "'"
$wsOverview.Cells.Item(3, 4).Style = "Normal"
$wsOverview.Cells.Item(3, 5).Value = $readyStatus
$wsOverview.Cells.Item(3, 6).Value = $readyStatus
$wsOverview.Cells.Item(3, 7).Value = $handoffDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $newHoUrl, "", "", $newPathName) | Out-Null
$wsOverview.Range("B3").Style = "HyperLink"
$wsOverview.Range("G3").NumberFormat = $dateFmt

$wsOverview.Columns.Item(5).ColumnWidth = 16.25
$wsOverview.Columns.Item(6).ColumnWidth = 16.25

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Cells.Item(3, 1).Value = $newBaseName
$wsZhCn.Cells.Item(3, 2).Value = ".md"
$wsZhCn.Cells.Item(3, 3).Value = $readyStatus
$wsZhCn.Cells.Item(3, 4).Value = "e2e"
$wsZhCn.Cells.Item(3, 5).Value = "ht"
$wsZhCn.Cells.Item(3, 6).Value = "'False"
$wsZhCn.Cells.Item(3, 6).Style = "Normal"
$wsZhCn.Cells.Item(3, 7).Value = $zhXliffName
$wsZhCn.Cells.Item(3, 8).Value = $xliffDate
$wsZhCn.Cells.Item(3, 9).Value = "'"
$wsZhCn.Cells.Item(3, 9).Style = "Normal"
$wsZhCn.Cells.Item(3, 10).Value = "'"
$wsZhCn.Cells.Item(3, 10).Style = "Normal"
$wsZhCn.Cells.Item(3, 11).Value = $epoch
$wsZhCn.Cells.Item(3, 12).Value = "'"
$wsZhCn.Cells.Item(3, 12).Style = "Normal"
$wsZhCn.Cells.Item(3, 13).Value = "'True"
$wsZhCn.Cells.Item(3, 13).Style = "Normal"
$wsZhCn.Cells.Item(3, 14).Value = "'"
$wsZhCn.Cells.Item(3, 14).Style = "Normal"
$wsZhCn.Cells.Item(3, 15).Value = "'False"
$wsZhCn.Cells.Item(3, 15).Style = "Normal"
$wsZhCn.Cells.Item(3, 16).Value = "'"
$wsZhCn.Cells.Item(3, 16).Style = "Normal"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $newHoUrl, "", "", $newBaseName) | Out-Null
$wsZhCn.Range("A3").Style = "HyperLink"
$wsZhCn.Range("H3").NumberFormat = $dateFmt
$wsZhCn.Range("K3").NumberFormat = $dateFmt

$wsZhCn.Columns.Item(3).ColumnWidth = 16.25

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Cells.Item(3, 1).Value = $newBaseName
$wsDeDe.Cells.Item(3, 2).Value = ".md"
$wsDeDe.Cells.Item(3, 3).Value = $readyStatus
$wsDeDe.Cells.Item(3, 4).Value = "e2e"
$wsDeDe.Cells.Item(3, 5).Value = "ht"
$wsDeDe.Cells.Item(3, 6).Value = "'False"
$wsDeDe.Cells.Item(3, 6).Style = "Normal"
$wsDeDe.Cells.Item(3, 7).Value = $deXliffName
$wsDeDe.Cells.Item(3, 8).Value = $handoffDate
$wsDeDe.Cells.Item(3, 9).Value = "'"
$wsDeDe.Cells.Item(3, 9).Style = "Normal"
$wsDeDe.Cells.Item(3, 10).Value = "'"
$wsDeDe.Cells.Item(3, 10).Style = "Normal"
$wsDeDe.Cells.Item(3, 11).Value = $epoch
$wsDeDe.Cells.Item(3, 12).Value = "'"
$wsDeDe.Cells.Item(3, 12).Style = "Normal"
$wsDeDe.Cells.Item(3, 13).Value = "'True"
$wsDeDe.Cells.Item(3, 13).Style = "Normal"
$wsDeDe.Cells.Item(3, 14).Value = "'"
$wsDeDe.Cells.Item(3, 14).Style = "Normal"
$wsDeDe.Cells.Item(3, 15).Value = "'False"
$wsDeDe.Cells.Item(3, 15).Style = "Normal"
$wsDeDe.Cells.Item(3, 16).Value = "'"
$wsDeDe.Cells.Item(3, 16).Style = "Normal"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $newHoUrl, "", "", $newBaseName) | Out-Null
$wsDeDe.Range("A3").Style = "HyperLink"
$wsDeDe.Range("H3").NumberFormat = $dateFmt
$wsDeDe.Range("K3").NumberFormat = $dateFmt

$wsDeDe.Columns.Item(3).ColumnWidth = 16.25

Write-Host "Handoff report rows appended."
